$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("day 1")

# Row 101 (S.No = 100): finish off the previously-partial row with the
# rest of the day's log entry.
$ws.Range("A101").Value = 100
$ws.Range("B101").Value = "8th Mar,2017"
$ws.Range("C101").Value = "Testing webflow till registration page"
$ws.Range("D101").Value = "NA"
$ws.Range("E101").Value = "NA"
$ws.Range("F101").Value = "120 minutes"
$ws.Range("G101").Value = "N"
$ws.Range("H101").Value = "NA"

# Row 102 (S.No = 101)
$ws.Range("A102").Value = 101
$ws.Range("B102").Value = "8th Mar,2017"
$ws.Range("C102").Value = "Created Address DAO,DTO and DAOIMPL and tested"
$ws.Range("D102").Value = "NA"
$ws.Range("E102").Value = "NA"
$ws.Range("F102").Value = "60 minutes"
$ws.Range("G102").Value = "N"
$ws.Range("H102").Value = "NA"

# Row 103 (S.No = 102)
$ws.Range("A103").Value = 102
$ws.Range("B103").Value = "8th Mar,2017"
$ws.Range("C103").Value = "Configured billing page with webflow"
$ws.Range("D103").Value = "NA"
$ws.Range("E103").Value = "NA"
$ws.Range("F103").Value = "60 minutes"
$ws.Range("G103").Value = "N"
$ws.Range("H103").Value = "NA"

# These three rows wrap onto two lines, matching the other multi-line rows
# in this log (e.g. rows 4-7, 9, 11-13, ...).
$ws.Range("A101:H103").EntireRow.RowHeight = 28.8

# Keep the sheet's dimension/selection in sync with the newly-added rows,
# same as Excel does after typing into new cells at the bottom of the
# used range.
[void]$ws.Range("A103").Select()
